$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("G:G").Insert()
$ws.Range("G5").Value = "{emailBenachrichtigungKiBonMail}"
$ws.Range("G4").Value = "{emailBenachrichtigungKiBonMailTitle}"
